# Update the answer values in the division-facts table.
# Each data row (1, 5, 9, 13, 17) holds 5 answer cells; we target each
# cell directly by (row, column) to avoid any ambiguity from duplicate
# text appearing elsewhere in the document (Find/Replace could otherwise
# cross-match a newly written value with an older, not-yet-updated cell).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "58÷7=8, 2"
$t.Cell(1, 2).Range.Text  = "45÷4=11, 1"
$t.Cell(1, 3).Range.Text  = "39÷3=13, 0"
$t.Cell(1, 4).Range.Text  = "57÷4=14, 1"
$t.Cell(1, 5).Range.Text  = "78÷7=11, 1"

$t.Cell(5, 1).Range.Text  = "39÷9=4, 3"
$t.Cell(5, 2).Range.Text  = "10÷6=1, 4"
$t.Cell(5, 3).Range.Text  = "59÷8=7, 3"
$t.Cell(5, 4).Range.Text  = "80÷8=10, 0"
$t.Cell(5, 5).Range.Text  = "30÷3=10, 0"

$t.Cell(9, 1).Range.Text  = "72÷4=18, 0"
$t.Cell(9, 2).Range.Text  = "99÷5=19, 4"
$t.Cell(9, 3).Range.Text  = "91÷9=10, 1"
$t.Cell(9, 4).Range.Text  = "97÷3=32, 1"
$t.Cell(9, 5).Range.Text  = "32÷9=3, 5"

$t.Cell(13, 1).Range.Text = "55÷5=11, 0"
$t.Cell(13, 2).Range.Text = "78÷8=9, 6"
$t.Cell(13, 3).Range.Text = "46÷4=11, 2"
$t.Cell(13, 4).Range.Text = "55÷4=13, 3"
$t.Cell(13, 5).Range.Text = "12÷8=1, 4"

$t.Cell(17, 1).Range.Text = "64÷4=16, 0"
$t.Cell(17, 2).Range.Text = "71÷6=11, 5"
$t.Cell(17, 3).Range.Text = "14÷6=2, 2"
$t.Cell(17, 4).Range.Text = "37÷8=4, 5"
$t.Cell(17, 5).Range.Text = "95÷4=23, 3"
